$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.041701534829099
$ws.Range("D2").Value = 1.051807231337337
$ws.Range("E2").Value = 1.050689735578913
$ws.Range("F2").Value = 1.063351488991097
$ws.Range("I2").Value = 1.046884152362405
$ws.Range("J2").Value = 1.046781323366449
$ws.Range("K2").Value = 1.05455771043539
$ws.Range("L2").Value = 1.053443314055529
$ws.Range("M2").Value = 1.066070362796151
$ws.Range("N2").Value = 1.019501132096821
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.04249508081187
$ws.Range("D3").Value = 1.052445151042215
$ws.Range("E3").Value = 1.051387029630809
$ws.Range("F3").Value = 1.064109725067819
$ws.Range("I3").Value = 1.047097409653571
$ws.Range("J3").Value = 1.04722178243274
$ws.Range("K3").Value = 1.055008799420326
$ws.Range("L3").Value = 1.053953401648132
$ws.Range("M3").Value = 1.066643735282242
$ws.Range("N3").Value = 1.019648030252441
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.043009244675728
$ws.Range("D4").Value = 1.052858494436539
$ws.Range("E4").Value = 1.051839193833649
$ws.Range("F4").Value = 1.064601333083946
$ws.Range("I4").Value = 1.047234450090868
$ws.Range("J4").Value = 1.047506780071183
$ws.Range("K4").Value = 1.055300560231129
$ws.Range("L4").Value = 1.054283751888188
$ws.Range("M4").Value = 1.067015056581838
$ws.Range("N4").Value = 1.019743055842762
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.043225561643893
$ws.Range("D5").Value = 1.05303239770519
$ws.Range("E5").Value = 1.052029513471146
$ws.Range("F5").Value = 1.064808237218747
$ws.Range("I5").Value = 1.047291833236413
$ws.Range("J5").Value = 1.047626589604162
$ws.Range("K5").Value = 1.055423185336707
$ws.Range("L5").Value = 1.054422698876066
$ws.Range("M5").Value = 1.067171232755397
$ws.Range("N5").Value = 1.019782997566594
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.043261891687908
$ws.Range("D6").Value = 1.053061604612707
$ws.Range("E6").Value = 1.052061482393137
$ws.Range("F6").Value = 1.064842990882829
$ws.Range("I6").Value = 1.047301454684115
$ws.Range("J6").Value = 1.047646705908644
$ws.Range("K6").Value = 1.055443772762696
$ws.Range("L6").Value = 1.054446032614343
$ws.Range("M6").Value = 1.067197459634484
$ws.Range("N6").Value = 1.019789703531044
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.043012134477565
$ws.Range("D7").Value = 1.052860817617622
$ws.Range("E7").Value = 1.05184173599351
$ws.Range("F7").Value = 1.064604096837731
$ws.Range("I7").Value = 1.047235217746862
$ws.Range("J7").Value = 1.047508380987386
$ws.Range("K7").Value = 1.055302198877994
$ws.Range("L7").Value = 1.054285608240389
$ws.Range("M7").Value = 1.067017143132049
$ws.Range("N7").Value = 1.019743589574285
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.041969574122787
$ws.Range("D8").Value = 1.052022700773578
$ws.Range("E8").Value = 1.050925187661419
$ws.Range("F8").Value = 1.063607534512839
$ws.Range("I8").Value = 1.046956420083795
$ws.Range("J8").Value = 1.046930179580791
$ws.Range("K8").Value = 1.054710182845255
$ws.Range("L8").Value = 1.053615639679671
$ws.Range("M8").Value = 1.066264071276865
$ws.Range("N8").Value = 1.019550782383017
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040137791027427
$ws.Range("D9").Value = 1.050550258630301
$ws.Range("E9").Value = 1.049317617554
$ws.Range("F9").Value = 1.061859054014193
$ws.Range("I9").Value = 1.046457897680662
$ws.Range("J9").Value = 1.045911307354581
$ws.Range("K9").Value = 1.053666090884811
$ws.Range("L9").Value = 1.052437354921933
$ws.Range("M9").Value = 1.064939521214469
$ws.Range("N9").Value = 1.019210843718926
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038920306042679
$ws.Range("D10").Value = 1.049571724045355
$ws.Range("E10").Value = 1.048251070488371
$ws.Range("F10").Value = 1.060698631588554
$ws.Range("I10").Value = 1.046120732720553
$ws.Range("J10").Value = 1.045232135453146
$ws.Range("K10").Value = 1.052969523087623
$ws.Range("L10").Value = 1.051653460994522
$ws.Range("M10").Value = 1.064058245299311
$ws.Range("N10").Value = 1.018984120926386
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.03839402384779
$ws.Range("D11").Value = 1.049148766335227
$ws.Range("E11").Value = 1.047790494727386
$ws.Range("F11").Value = 1.060197423863214
$ws.Range("I11").Value = 1.045973605658849
$ws.Range("J11").Value = 1.044938081803183
$ws.Range("K11").Value = 1.052667800013285
$ws.Range("L11").Value = 1.051314432324949
$ws.Range("M11").Value = 1.06367708204159
$ws.Range("N11").Value = 1.018885930627315
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.038198675869819
$ws.Range("D12").Value = 1.04899177634282
$ws.Range("E12").Value = 1.047619605684348
$ws.Range("F12").Value = 1.060011444943596
$ws.Range("I12").Value = 1.045918786904871
$ws.Range("J12").Value = 1.04482886342239
$ws.Range("K12").Value = 1.052555712414563
$ws.Range("L12").Value = 1.051188564056284
$ws.Range("M12").Value = 1.063535568187763
$ws.Range("N12").Value = 1.018849456199141
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.038240572459099
$ws.Range("D13").Value = 1.049025445977634
$ws.Range("E13").Value = 1.047656253364501
$ws.Range("F13").Value = 1.060051329331639
$ws.Range("I13").Value = 1.045930553356134
$ws.Range("J13").Value = 1.044852290831249
$ws.Range("K13").Value = 1.052579756204096
$ws.Range("L13").Value = 1.051215560392636
$ws.Range("M13").Value = 1.063565920322276
$ws.Range("N13").Value = 1.018857280178388
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.038377873524007
$ws.Range("D14").Value = 1.049135787129047
$ws.Range("E14").Value = 1.047776365112061
$ws.Range("F14").Value = 1.060182046859248
$ws.Range("I14").Value = 1.045969077770106
$ws.Range("J14").Value = 1.044929053642476
$ws.Range("K14").Value = 1.052658535098694
$ws.Range("L14").Value = 1.051304026739094
$ws.Range("M14").Value = 1.063665383082055
$ws.Range("N14").Value = 1.018882915681268
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.03846248740816
$ws.Range("D15").Value = 1.049203787297294
$ws.Range("E15").Value = 1.047850395070415
$ws.Range("F15").Value = 1.060262611743212
$ws.Range("I15").Value = 1.045992791538249
$ws.Range("J15").Value = 1.044976350612064
$ws.Range("K15").Value = 1.052707071537054
$ws.Range("L15").Value = 1.05135854204797
$ws.Range("M15").Value = 1.063726674322855
$ws.Range("N15").Value = 1.018898710286858
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.038955252718065
$ws.Range("D16").Value = 1.049599810418696
$ws.Range("E16").Value = 1.048281663817537
$ws.Range("F16").Value = 1.060731921890634
$ws.Range("I16").Value = 1.046130473278966
$ws.Range("J16").Value = 1.045251651623916
$ws.Range("K16").Value = 1.052989545378475
$ws.Range("L16").Value = 1.051675969828484
$ws.Range("M16").Value = 1.064083551168065
$ws.Range("N16").Value = 1.018990637160123
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039264592843603
$ws.Range("D17").Value = 1.049848428613696
$ws.Range("E17").Value = 1.04855252241008
$ws.Range("F17").Value = 1.061026647412581
$ws.Range("I17").Value = 1.046216534867865
$ws.Range("J17").Value = 1.045424350295138
$ws.Range("K17").Value = 1.05316670680789
$ws.Range("L17").Value = 1.051875192793918
$ws.Range("M17").Value = 1.064307528320566
$ws.Range("N17").Value = 1.019048296037705
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039445112108043
$ws.Range("D18").Value = 1.049993516036805
$ws.Range("E18").Value = 1.048710629741223
$ws.Range("F18").Value = 1.061198677520576
$ws.Range("I18").Value = 1.046266623813858
$ws.Range("J18").Value = 1.045525085528178
$ws.Range("K18").Value = 1.053270031877104
$ws.Range("L18").Value = 1.051991434914924
$ws.Range("M18").Value = 1.064438212258067
$ws.Range("N18").Value = 1.019081925731281
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039506679109147
$ws.Range("D19").Value = 1.050042999357558
$ws.Range("E19").Value = 1.04876456056241
$ws.Range("F19").Value = 1.061257355916217
$ws.Range("I19").Value = 1.046283684279043
$ws.Range("J19").Value = 1.045559434109348
$ws.Range("K19").Value = 1.053305261280503
$ws.Range("L19").Value = 1.052031076996028
$ws.Range("M19").Value = 1.064482779148135
$ws.Range("N19").Value = 1.019093392274953
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039231394624788
$ws.Range("D20").Value = 1.049821746706596
$ws.Range("E20").Value = 1.048523449409313
$ws.Range("F20").Value = 1.060995013549604
$ws.Range("I20").Value = 1.046207312574996
$ws.Range("J20").Value = 1.045405821039729
$ws.Range("K20").Value = 1.053147700106851
$ws.Range("L20").Value = 1.051853814037041
$ws.Range("M20").Value = 1.064283493349732
$ws.Range("N20").Value = 1.019042109963658
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.038337437994547
$ws.Range("D21").Value = 1.049103291210167
$ws.Range("E21").Value = 1.047740989953763
$ws.Range("F21").Value = 1.060143548488973
$ws.Range("I21").Value = 1.045957737954089
$ws.Range("J21").Value = 1.044906448724486
$ws.Range("K21").Value = 1.052635337058073
$ws.Range("L21").Value = 1.05127797387869
$ws.Range("M21").Value = 1.063636091910245
$ws.Range("N21").Value = 1.018875366720761
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.037776163616579
$ws.Range("D22").Value = 1.048652238004649
$ws.Range("E22").Value = 1.04725012290726
$ws.Range("F22").Value = 1.059609310246098
$ws.Range("I22").Value = 1.045799841929356
$ws.Range("J22").Value = 1.044592510305861
$ws.Range("K22").Value = 1.052313112575298
$ws.Range("L22").Value = 1.050916280246731
$ws.Range("M22").Value = 1.063229434333454
$ws.Range("N22").Value = 1.018770516265005
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.038073630103705
$ws.Range("D23").Value = 1.048891285824225
$ws.Range("E23").Value = 1.047510236208581
$ws.Range("F23").Value = 1.059892413936711
$ws.Range("I23").Value = 1.045883638047984
$ws.Range("J23").Value = 1.044758931051593
$ws.Range("K23").Value = 1.05248393709365
$ws.Range("L23").Value = 1.051107986236008
$ws.Range("M23").Value = 1.063444973699883
$ws.Range("N23").Value = 1.01882610048027
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039246395204148
$ws.Range("D24").Value = 1.049833802887743
$ws.Range("E24").Value = 1.048536585876129
$ws.Range("F24").Value = 1.061009307154893
$ws.Range("I24").Value = 1.046211480070081
$ws.Range("J24").Value = 1.045414193603517
$ws.Range("K24").Value = 1.053156288448616
$ws.Range("L24").Value = 1.051863474057559
$ws.Range("M24").Value = 1.064294353589395
$ws.Range("N24").Value = 1.019044905189729
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040610706294539
$ws.Range("D25").Value = 1.050930383533744
$ws.Range("E25").Value = 1.049732311153892
$ws.Range("F25").Value = 1.062310165483217
$ws.Range("I25").Value = 1.046587630651475
$ws.Range("J25").Value = 1.046174702755954
$ws.Range("K25").Value = 1.053936109388732
$ws.Range("L25").Value = 1.052741689386825
$ws.Range("M25").Value = 1.065281647023143
$ws.Range("N25").Value = 1.019298745249171